# Project Sample Project is saved. The "Rules" decision table row 11 (R40)
# is re-keyed from the text label "R40" to the text label "1" in column B.
#
# Note: we must land a *text* value ("1") in B11 while leaving its style
# (General number format, no quote-prefix) untouched. A plain
# `Range.Value = "1"` assignment would be auto-coerced to the number 1
# (Excel's normal "typed into a General cell" behaviour), so instead we
# stage the text in a scratch cell via a formula that yields a string,
# copy it, and paste-special just the *values* into B11 - this carries
# the string type over without touching number formatting. The scratch
# cell is fully cleared afterwards so no trace of it remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("G1")
$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)
$scratch.Clear()
